$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-253)
# from serial date 45202 to 45203.
$ws.Range("C2:C253").Value = 45203
